$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")

# Add the "project surfix" values in column D for rows 119-138
$values = @{
    119 = "guacamole-server"
    120 = "hama"
    121 = "harmony"
    122 = "hbase"
    123 = "helix"
    124 = "hive"
    125 = "hivemind"
    126 = "hop"
    127 = "httpcomponents"
    128 = "httpcomponents"
    129 = "hudi"
    130 = "ignite"
    131 = "iotdb"
    132 = "isis"
    133 = "ivy"
    134 = "ivyde"
    135 = "jackrabbit"
    136 = "jakarta cactus"
    137 = "james"
    138 = "jclouds"
}

foreach ($row in 119..138) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# Update the sheet's active selection to match the edited workbook state
$ws.Activate()
$ws.Range("E131").Select()
